# Apply formatting changes to Sheet1 of the Customer workbook:
# - Column A (rows 1-22): solid fill theme color Accent1 (theme index 4), thin black border
# - Column B (rows 1-22): solid fill RGB FFC000 (orange), thin black border, incl. new empty B11
# - Update selection to C11:D11

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$colA = $ws.Range("A1:A22")
$colB = $ws.Range("B1:B22")
$colAB = $ws.Range("A1:B22")

# Apply the thin black border to both columns first
$colAB.Borders.LineStyle = 1     # xlContinuous
$colAB.Borders.Weight = 2        # xlThin

# Column A fill: theme color Accent1 (theme index 4)
$colA.Interior.ThemeColor = 5    # msoThemeColorAccent1 -> theme="4" in the xf

# Column B fill: solid orange RGB(255,192,0) = FFC000
$colB.Interior.Color = 49407

# Update the active selection on Sheet1 to C11:D11
$ws.Range("C11:D11").Select()
